$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Assignment cell (B7) from "Bài tập nhóm Tuần 1" to "Bài tập nhóm Tuần 2"
$ws.Range("B7").Value = "Bài tập nhóm Tuần 2"

# Widen column A slightly (22.85546875 -> ~25.71 chars)
$ws.Columns.Item(1).ColumnWidth = 24.8

# Move the active selection from F7 to E25
$ws.Range("E25").Select() | Out-Null
